$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell($table, $row, $col, $a, $b) {
    $aStr = [string]$a
    $bStr = [string]$b
    $d1 = $bStr.Substring(0,1)
    $d2 = $bStr.Substring(1,1)
    $r1 = $aStr.Substring(0,1)
    $r2 = $aStr.Substring(1,1)

    $text = "$aStr x $bStr" + "`v" + "  $d1    $d2" + "`v" + "  ----" + "`v" + "$r1|    |" + "`v" + "$r2|    |"

    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

Set-LatticeCell $t 1 1 97 62
Set-LatticeCell $t 1 2 72 75
Set-LatticeCell $t 1 3 11 88

Set-LatticeCell $t 2 1 39 76
Set-LatticeCell $t 2 2 61 40
Set-LatticeCell $t 2 3 28 69

Set-LatticeCell $t 3 1 98 11
Set-LatticeCell $t 3 2 35 21
Set-LatticeCell $t 3 3 70 75

Set-LatticeCell $t 4 1 52 34
Set-LatticeCell $t 4 2 36 12
Set-LatticeCell $t 4 3 55 67

Set-LatticeCell $t 5 1 81 14
Set-LatticeCell $t 5 2 95 36
Set-LatticeCell $t 5 3 30 79
